$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (student #15): fill in the previously-blank marks with 5
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 5

# Row 29 (student #26): fill in the remaining blank marks with 5
# (C29 already had a 5, D29/E29/F29 were blank)
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 5

# Move the active selection in the (still frozen at C4 / xSplit=2,ySplit=3)
# bottom-right pane to F29, scrolling the frozen pane's viewport down so
# row 16 is the first visible row (topLeftCell C16).
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 3
$ws.Range("F29").Select()
